$wb = $excel.ActiveWorkbook

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 146.55556
$ws.Cells.Item(33, 9).Value = 103.166664
$ws.Cells.Item(33, 11).Value = 103.166664
$ws.Cells.Item(33, 13).Value = 125.833336

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 1080
$ws.Cells.Item(98, 9).Value = 1038.5
$ws.Cells.Item(98, 11).Value = 1038.5
$ws.Cells.Item(98, 13).Value = 459.5

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 4999.5
$ws.Cells.Item(113, 9).Value = 5000
$ws.Cells.Item(113, 10).Value = 4999
$ws.Cells.Item(113, 11).Value = 5000
$ws.Cells.Item(113, 12).Value = 4999
$ws.Cells.Item(113, 13).Value = -1746
$ws.Cells.Item(113, 14).Value = -11507

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 1080
$ws.Cells.Item(122, 9).Value = 1038.5
$ws.Cells.Item(122, 11).Value = 3115.5
$ws.Cells.Item(122, 13).Value = -665.5

# ARM row 15
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).ClearContents()
$ws.Cells.Item(15, 13).ClearContents()
$ws.Cells.Item(15, 14).Value = 0

# ARM row 17
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).ClearContents()
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(17, 14).Value = 0

# ARM row 58
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 12).ClearContents()
$ws.Cells.Item(58, 14).Value = 0

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1811.5385
$ws.Cells.Item(61, 9).Value = 1825
$ws.Cells.Item(61, 11).Value = 1825
$ws.Cells.Item(61, 13).Value = -1613

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 2499.3076
$ws.Cells.Item(110, 9).Value = 1938
$ws.Cells.Item(110, 10).Value = 3154.1667
$ws.Cells.Item(110, 11).Value = 1938
$ws.Cells.Item(110, 12).Value = 3154.1667
$ws.Cells.Item(110, 13).Value = 107
$ws.Cells.Item(110, 14).Value = -7244.1667

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2969
$ws.Cells.Item(122, 9).Value = 2914.6667
$ws.Cells.Item(122, 11).Value = 8744.000100000001
$ws.Cells.Item(122, 13).Value = -6294.000100000001

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1811.5385
$ws.Cells.Item(136, 9).Value = 1825
$ws.Cells.Item(136, 11).Value = 5475
$ws.Cells.Item(136, 13).Value = -2925

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2301.3333
$ws.Cells.Item(20, 9).Value = 2299.8
$ws.Cells.Item(20, 11).Value = 2299.8
$ws.Cells.Item(20, 13).Value = -2052.8

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 585.4286
$ws.Cells.Item(94, 9).Value = 499.6
$ws.Cells.Item(94, 11).Value = 499.6
$ws.Cells.Item(94, 13).Value = -48.60000000000002

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2068.3333
$ws.Cells.Item(99, 9).Value = 2089.375
$ws.Cells.Item(99, 10).Value = 1900
$ws.Cells.Item(99, 11).Value = 2089.375
$ws.Cells.Item(99, 12).Value = 1900
$ws.Cells.Item(99, 13).Value = -591.375
$ws.Cells.Item(99, 14).Value = -4896

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 5438.8667
$ws.Cells.Item(107, 9).Value = 1297.5714
$ws.Cells.Item(107, 11).Value = 1297.5714
$ws.Cells.Item(107, 13).Value = 622.4286

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2304.125
$ws.Cells.Item(134, 9).Value = 2304.125
$ws.Cells.Item(134, 11).Value = 6912.375
$ws.Cells.Item(134, 13).Value = -4377.375

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 3195.8
$ws.Cells.Item(22, 9).Value = 1114.5
$ws.Cells.Item(22, 10).Value = 4583.3335
$ws.Cells.Item(22, 11).Value = 1114.5
$ws.Cells.Item(22, 12).Value = 4583.3335
$ws.Cells.Item(22, 13).Value = -764.5
$ws.Cells.Item(22, 14).Value = -5283.3335

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6334.032
$ws.Cells.Item(31, 9).Value = 4429.8335
$ws.Cells.Item(31, 11).Value = 4429.8335
$ws.Cells.Item(31, 13).Value = -4134.8335

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 6334.032
$ws.Cells.Item(34, 9).Value = 4429.8335
$ws.Cells.Item(34, 11).Value = 4429.8335
$ws.Cells.Item(34, 13).Value = -4227.8335

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 3496.4
$ws.Cells.Item(58, 9).Value = 1710.5714
$ws.Cells.Item(58, 10).Value = 7663.3335
$ws.Cells.Item(58, 11).Value = 1710.5714
$ws.Cells.Item(58, 12).Value = 7663.3335
$ws.Cells.Item(58, 13).Value = -1507.5714
$ws.Cells.Item(58, 14).Value = -8069.3335

# CRP row 82
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(82, 8).Value = 6000
$ws.Cells.Item(82, 9).Value = 6000
$ws.Cells.Item(82, 11).Value = 6000
$ws.Cells.Item(82, 13).Value = -5639

# CRP row 85
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(85, 8).Value = 6000
$ws.Cells.Item(85, 9).Value = 6000
$ws.Cells.Item(85, 11).Value = 6000
$ws.Cells.Item(85, 13).Value = -4752

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(95, 8).Value = 13333.333
$ws.Cells.Item(95, 10).Value = 13333.333
$ws.Cells.Item(95, 12).Value = 13333.333
$ws.Cells.Item(95, 14).Value = -18825.333

# CRP row 96
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(96, 8).Value = 7327.625
$ws.Cells.Item(96, 10).Value = 7327.625
$ws.Cells.Item(96, 12).Value = 7327.625
$ws.Cells.Item(96, 14).Value = -12819.625

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2978.5
$ws.Cells.Item(132, 9).Value = 2842.7778
$ws.Cells.Item(132, 11).Value = 8528.3334
$ws.Cells.Item(132, 13).Value = -5998.3334

# CRP row 133
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(133, 8).Value = 45000
$ws.Cells.Item(133, 10).Value = 45000
$ws.Cells.Item(133, 12).Value = 45000
$ws.Cells.Item(133, 14).Value = -50060

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2980.3333
$ws.Cells.Item(134, 9).Value = 2980.3333
$ws.Cells.Item(134, 11).Value = 8940.999899999999
$ws.Cells.Item(134, 13).Value = -6405.999899999999

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 3496.4
$ws.Cells.Item(136, 9).Value = 1710.5714
$ws.Cells.Item(136, 10).Value = 7663.3335
$ws.Cells.Item(136, 11).Value = 5131.7142
$ws.Cells.Item(136, 12).Value = 22990.0005
$ws.Cells.Item(136, 13).Value = -2581.7142
$ws.Cells.Item(136, 14).Value = -28090.0005

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(140, 8).Value = 148947.6
$ws.Cells.Item(140, 10).Value = 148947.6
$ws.Cells.Item(140, 12).Value = 148947.6
$ws.Cells.Item(140, 14).Value = -159307.6

# CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 29.411764
$ws.Cells.Item(2, 9).Value = 39
$ws.Cells.Item(2, 10).Value = 20.88889
$ws.Cells.Item(2, 11).Value = 234
$ws.Cells.Item(2, 12).Value = 125.33334
$ws.Cells.Item(2, 13).Value = -121
$ws.Cells.Item(2, 14).Value = -351.33334

# CUL row 57
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(57, 9).Value = 2000
$ws.Cells.Item(57, 10).Value = 2750
$ws.Cells.Item(57, 11).Value = 6000
$ws.Cells.Item(57, 12).Value = 8250
$ws.Cells.Item(57, 13).Value = -5441
$ws.Cells.Item(57, 14).Value = -9368

# CUL row 60
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(60, 8).Value = 1068.4286
$ws.Cells.Item(60, 9).Value = 292.18182
$ws.Cells.Item(60, 10).Value = 1922.3
$ws.Cells.Item(60, 11).Value = 876.54546
$ws.Cells.Item(60, 12).Value = 5766.9
$ws.Cells.Item(60, 13).Value = -625.54546
$ws.Cells.Item(60, 14).Value = -6268.9

# CUL row 128
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(128, 8).Value = 508328.34
$ws.Cells.Item(128, 9).Value = 508328.34
$ws.Cells.Item(128, 11).Value = 1524985.02
$ws.Cells.Item(128, 13).Value = -1520005.02

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(139, 8).Value = 1875.6154
$ws.Cells.Item(139, 9).Value = 1130.8
$ws.Cells.Item(139, 11).Value = 3392.4
$ws.Cells.Item(139, 13).Value = 1747.6

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2233
$ws.Cells.Item(80, 9).Value = 2233
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 11).Value = 2233
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
$ws.Cells.Item(80, 14).Value = -1235

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 2233
$ws.Cells.Item(83, 9).Value = 2233
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 11).Value = 11165
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
$ws.Cells.Item(83, 14).Value = -6173

# GSM row 101
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(101, 8).Value = 26666
$ws.Cells.Item(101, 10).Value = 26666
$ws.Cells.Item(101, 12).Value = 26666
$ws.Cells.Item(101, 14).Value = -33156

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1704
$ws.Cells.Item(102, 9).Value = 1704
$ws.Cells.Item(102, 11).Value = 1704
$ws.Cells.Item(102, 13).Value = -82

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2690.5
$ws.Cells.Item(126, 9).Value = 1535.75
$ws.Cells.Item(126, 11).Value = 4607.25
$ws.Cells.Item(126, 13).Value = -2137.25

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 96385.09
$ws.Cells.Item(132, 9).Value = 204307.6
$ws.Cells.Item(132, 10).Value = 6449.6665
$ws.Cells.Item(132, 11).Value = 612922.8
$ws.Cells.Item(132, 12).Value = 19348.9995
$ws.Cells.Item(132, 13).Value = -610392.8
$ws.Cells.Item(132, 14).Value = -24408.9995

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 5824.75
$ws.Cells.Item(61, 10).Value = 6333
$ws.Cells.Item(61, 12).Value = 6333
$ws.Cells.Item(61, 14).Value = -6737

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 5645.364
$ws.Cells.Item(100, 9).Value = 4012.375
$ws.Cells.Item(100, 10).Value = 10000
$ws.Cells.Item(100, 11).Value = 4012.375
$ws.Cells.Item(100, 12).Value = 10000
$ws.Cells.Item(100, 13).Value = -3471.375
$ws.Cells.Item(100, 14).Value = -11082

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 5824.75
$ws.Cells.Item(113, 10).Value = 6333
$ws.Cells.Item(113, 12).Value = 6333
$ws.Cells.Item(113, 14).Value = -10673

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2849.75
$ws.Cells.Item(122, 9).Value = 2299.6667
$ws.Cells.Item(122, 11).Value = 6899.000100000001
$ws.Cells.Item(122, 13).Value = -4449.000100000001

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 2428.9
$ws.Cells.Item(136, 9).Value = 2428.9
$ws.Cells.Item(136, 11).Value = 7286.700000000001
$ws.Cells.Item(136, 13).Value = -4736.700000000001

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8505.454
$ws.Cells.Item(62, 9).Value = 5976
$ws.Cells.Item(62, 10).Value = 9454
$ws.Cells.Item(62, 11).Value = 5976
$ws.Cells.Item(62, 12).Value = 9454
$ws.Cells.Item(62, 13).Value = -5352
$ws.Cells.Item(62, 14).Value = -10702

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 8505.454
$ws.Cells.Item(65, 9).Value = 5976
$ws.Cells.Item(65, 10).Value = 9454
$ws.Cells.Item(65, 11).Value = 29880
$ws.Cells.Item(65, 12).Value = 47270
$ws.Cells.Item(65, 13).Value = -26760
$ws.Cells.Item(65, 14).Value = -53510

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 678.5714
$ws.Cells.Item(81, 9).Value = 625
$ws.Cells.Item(81, 11).Value = 1250
$ws.Cells.Item(81, 13).Value = -189

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 678.5714
$ws.Cells.Item(84, 9).Value = 625
$ws.Cells.Item(84, 11).Value = 6250
$ws.Cells.Item(84, 13).Value = -946

# WVR row 95
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 17814.666
$ws.Cells.Item(95, 10).Value = 17814.666
$ws.Cells.Item(95, 12).Value = 17814.666
$ws.Cells.Item(95, 14).Value = -23306.666

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1644.5
$ws.Cells.Item(122, 9).Value = 1644.5
$ws.Cells.Item(122, 11).Value = 4933.5
$ws.Cells.Item(122, 13).Value = -2483.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 4000
$ws.Cells.Item(132, 9).Value = 3500
$ws.Cells.Item(132, 11).Value = 10500
$ws.Cells.Item(132, 13).Value = -7970
